$d = $word.ActiveDocument

# The target paragraph currently holds its text split across four runs:
#   "Change the date to th" + "e" + " end of march" + " (PRIO 3)"
# Collapse them into a single run with the full sentence by doing a
# literal find & replace over the exact concatenated text.
$findText = "Change the date to th" + "e" + " end of march" + " (PRIO 3)"

$rng = $d.Content
$rng.Find.Execute($findText, $true, $false, $false, $false, $false, `
                   $true, 1, $false, $findText, 2)

# Now find that paragraph and apply a green highlight to its whole range
# (this covers every run in the paragraph as well as the paragraph mark,
# matching how Word stamps <w:rPr> on both the runs and the enclosing
# <w:pPr> when a full paragraph selection is highlighted).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Change the date to the end of march*") {
        $p.Range.Font.HighlightColorIndex = 4
    }
}
